# Apply the target edit: fill previously-empty numeric cells on two worksheets.
$wb = $excel.ActiveWorkbook

# --- "Hyperparameter Optimization" sheet: columns S,T,U,V for rows 72,73,75,76,77,79,80,81,83 ---
$wsHyper = $wb.Worksheets.Item("Hyperparameter Optimization")

$wsHyper.Range("S72").Value = 0.577744205979862
$wsHyper.Range("T72").Value = 0.465474121633776
$wsHyper.Range("U72").Value = 0.104045624312589
$wsHyper.Range("V72").Value = 0.905099999518191
$wsHyper.Range("S73").Value = 0.946146326111374
$wsHyper.Range("T73").Value = 0.717107629990938
$wsHyper.Range("U73").Value = 0.151137137162642
$wsHyper.Range("V73").Value = 0.748157560711823
$wsHyper.Range("S75").Value = 0.9088329510916
$wsHyper.Range("T75").Value = 0.680121840787656
$wsHyper.Range("U75").Value = 0.134257289961237
$wsHyper.Range("V75").Value = 0.845351991027389
$wsHyper.Range("S76").Value = 1.55205446517014
$wsHyper.Range("T76").Value = 1.09309311232192
$wsHyper.Range("U76").Value = 0.152265572776674
$wsHyper.Range("V76").Value = 0.858956934632023
$wsHyper.Range("S77").Value = 0.740865287597644
$wsHyper.Range("T77").Value = 0.537966951174266
$wsHyper.Range("U77").Value = 0.0840829448007078
$wsHyper.Range("V77").Value = 0.963147900757118
$wsHyper.Range("S79").Value = 0.72134524999995
$wsHyper.Range("T79").Value = 0.539329912489562
$wsHyper.Range("U79").Value = 0.0792728671958818
$wsHyper.Range("V79").Value = 0.963215368282948
$wsHyper.Range("S80").Value = 1.2266217502358
$wsHyper.Range("T80").Value = 0.971449108326706
$wsHyper.Range("U80").Value = 0.125268482682261
$wsHyper.Range("V80").Value = 0.925896624843258
$wsHyper.Range("S81").Value = 0.925071830869918
$wsHyper.Range("T81").Value = 0.697102822119073
$wsHyper.Range("U81").Value = 0.0912049940508836
$wsHyper.Range("V81").Value = 0.961490511073966
$wsHyper.Range("S83").Value = 0.902294505410523
$wsHyper.Range("T83").Value = 0.72373146750205
$wsHyper.Range("U83").Value = 0.0966128284047383
$wsHyper.Range("V83").Value = 0.960631029131909

# --- "CNNLSTM SH" sheet: columns B,C,D,E,F and I for rows 10,11,13,16,17,19,22,23,25 ---
# (G10:G25 and the summary rows 14/20/26 hold existing AVERAGE() formulas that
#  recompute automatically once B:F are populated, clearing the prior #DIV/0! errors.)
$wsCnn = $wb.Worksheets.Item("CNNLSTM SH")

$wsCnn.Range("B10").Value = 0.478084291556413
$wsCnn.Range("C10").Value = 0.545501079768484
$wsCnn.Range("D10").Value = 0.600558842303638
$wsCnn.Range("E10").Value = 0.493996689244862
$wsCnn.Range("F10").Value = 0.770580127025913
$wsCnn.Range("I10").Value = 0.11255421688019
$wsCnn.Range("B11").Value = 0.680072675822668
$wsCnn.Range("C11").Value = 0.647489141292674
$wsCnn.Range("D11").Value = 0.993761349121729
$wsCnn.Range("E11").Value = 1.07365172324007
$wsCnn.Range("F11").Value = 1.33575674107972
$wsCnn.Range("I11").Value = 0.184325100429041
$wsCnn.Range("B13").Value = 0.985651602753388
$wsCnn.Range("C13").Value = 0.688747416692525
$wsCnn.Range("D13").Value = 1.14542797153291
$wsCnn.Range("E13").Value = 0.676181531277256
$wsCnn.Range("F13").Value = 1.04815623320192
$wsCnn.Range("I13").Value = 0.177055832020915
$wsCnn.Range("B16").Value = 1.12713151536374
$wsCnn.Range("C16").Value = 1.73090106942843
$wsCnn.Range("D16").Value = 1.54073007232601
$wsCnn.Range("E16").Value = 1.89712801158729
$wsCnn.Range("F16").Value = 1.46438165714521
$wsCnn.Range("I16").Value = 0.20082103085223
$wsCnn.Range("B17").Value = 0.518586173021849
$wsCnn.Range("C17").Value = 0.750602942433554
$wsCnn.Range("D17").Value = 1.15329302134481
$wsCnn.Range("E17").Value = 0.703010785445734
$wsCnn.Range("F17").Value = 0.578833515742267
$wsCnn.Range("I17").Value = 0.0958608954239782
$wsCnn.Range("B19").Value = 0.461306619877098
$wsCnn.Range("C19").Value = 0.663236957233487
$wsCnn.Range("D19").Value = 1.31579189862773
$wsCnn.Range("E19").Value = 0.541290857728334
$wsCnn.Range("F19").Value = 0.625099916533096
$wsCnn.Range("I19").Value = 0.0933351889100554
$wsCnn.Range("B22").Value = 1.36071146058686
$wsCnn.Range("C22").Value = 0.950726209362045
$wsCnn.Range("D22").Value = 1.07039182408789
$wsCnn.Range("E22").Value = 1.70752953074928
$wsCnn.Range("F22").Value = 1.04374972639296
$wsCnn.Range("I22").Value = 0.131548259985608
$wsCnn.Range("B23").Value = 0.653753663161734
$wsCnn.Range("C23").Value = 0.92221970609107
$wsCnn.Range("D23").Value = 1.02698772201034
$wsCnn.Range("E23").Value = 1.11935725610742
$wsCnn.Range("F23").Value = 0.903040806979016
$wsCnn.Range("I23").Value = 0.0992087330012245
$wsCnn.Range("B25").Value = 0.464260937346463
$wsCnn.Range("C25").Value = 1.10544441771106
$wsCnn.Range("D25").Value = 0.85311994357041
$wsCnn.Range("E25").Value = 1.31546381861766
$wsCnn.Range("F25").Value = 0.773183409807011
$wsCnn.Range("I25").Value = 0.0967659933948762

# --- Restore the on-screen selections recorded in the target workbook ---
$wsCnn.Range("C23").Select()
$wsHyper.Activate()
$wsHyper.Range("S84").Select()
